$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.124.39'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.262.13'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'271.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.45%  '
$ws.Range("D6").Value = "'86.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +12.59%  '
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.44%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = "'0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.47%  '
$ws.Range("D10").Value = "'45.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.90%  '
$ws.Range("E11").Value = '  +1.80%  '
$ws.Range("D12").Value = "'7.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.05%  '
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").Value = '2.607.86'
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("D15").Value = "'15.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.69%  '
$ws.Range("D16").Value = '2.280.87'
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("D17").Value = "'0.804"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.23%  '
$ws.Range("D18").Value = '44.067.57'
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").Value = "'6.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = "'70.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").Value = "'235.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.20%  '
$ws.Range("D24").Value = "'8.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.54%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = "'2.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.99%  '
$ws.Range("D27").Value = "'10.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("E28").Value = '  +6.80%  '
$ws.Range("D29").Value = "'2.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.95%  '
$ws.Range("D30").Value = "'39.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.58%  '
$ws.Range("D31").Value = "'174.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = "'0.0904"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.50%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'20.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.90%  '
$ws.Range("D34").Value = "'5.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.68%  '
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("D36").Value = "'0.112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.07%  '
$ws.Range("D37").Value = "'0.0351"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.08%  '
$ws.Range("D38").Value = "'4.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("D39").Value = "'3.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +15.53%  '
$ws.Range("E40").Value = '  +5.77%  '
$ws.Range("D41").Value = "'12.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.73%  '
$ws.Range("D42").Value = "'64.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.76%  '
$ws.Range("E43").Value = '  +4.02%  '
$ws.Range("D44").Value = "'0.203"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = "'8.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'101.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.97%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = "'0.0990"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("E48").Value = '  +6.36%  '
$ws.Range("E49").Value = '  +2.68%  '
$ws.Range("D50").Value = "'1.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.35%  '
$ws.Range("D51").Value = "'0.432"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.25%  '
